# Update team transition-probability matrix (Westcliff_A) with values
# recomputed from games pulled March 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2666666666666667
$ws.Range("C2").Value = 0.4333333333333333
$ws.Range("P2").Value = 0.2333333333333333
$ws.Range("S2").Value = 0.06666666666666667

$ws.Range("J3").Value = 0.07692307692307693
$ws.Range("P3").Value = 0.6923076923076923
$ws.Range("S3").Value = 0.2307692307692308

$ws.Range("B6").Value = 0.1176470588235294
$ws.Range("F6").Value = 0.1176470588235294
$ws.Range("J6").Value = 0.2352941176470588
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.4705882352941176

$ws.Range("B7").Value = 0.1333333333333333
$ws.Range("J7").Value = 0.06666666666666667
$ws.Range("Q7").Value = 0.06666666666666667
$ws.Range("R7").Value = 0.06666666666666667
$ws.Range("S7").Value = 0.6666666666666666

$ws.Range("B8").Value = 0.06
$ws.Range("D8").Value = 0.02
$ws.Range("F8").Value = 0.06
$ws.Range("J8").Value = 0.1
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.06
$ws.Range("S8").Value = 0.5

$ws.Range("B9").Value = 0.1333333333333333
$ws.Range("E9").Value = 0.06666666666666667
$ws.Range("Q9").Value = 0.1333333333333333
$ws.Range("R9").Value = 0.1333333333333333
$ws.Range("S9").Value = 0.5333333333333333

$ws.Range("B10").Value = 0.1363636363636364
$ws.Range("D10").Value = 0.02272727272727273
$ws.Range("F10").Value = 0.06818181818181818
$ws.Range("J10").Value = 0.07954545454545454
$ws.Range("O10").Value = 0.03409090909090909
$ws.Range("Q10").Value = 0.1704545454545454
$ws.Range("R10").Value = 0.04545454545454546
$ws.Range("S10").Value = 0.4431818181818182

$ws.Range("G11").Value = 0.2173913043478261
$ws.Range("K11").Value = 0.2608695652173913
$ws.Range("L11").Value = 0.4782608695652174
$ws.Range("S11").Value = 0.04347826086956522

$ws.Range("G12").Value = 0.7272727272727273
$ws.Range("J12").Value = 0.09090909090909091
$ws.Range("S12").Value = 0.1818181818181818

$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.25

$ws.Range("H15").Value = 0.2857142857142857
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("O15").Value = 0.04761904761904762
$ws.Range("S15").Value = 0.2380952380952381

$ws.Range("F16").Value = 0.05263157894736842
$ws.Range("H16").Value = 0.2631578947368421
$ws.Range("I16").Value = 0.05263157894736842
$ws.Range("J16").Value = 0.3157894736842105
$ws.Range("K16").Value = 0.1578947368421053
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1052631578947368

$ws.Range("H17").Value = 0.2857142857142857
$ws.Range("I17").Value = 0.03571428571428571
$ws.Range("J17").Value = 0.4285714285714285
$ws.Range("K17").Value = 0.03571428571428571
$ws.Range("M17").Value = 0.03571428571428571
$ws.Range("O17").Value = 0.1071428571428571
$ws.Range("S17").Value = 0.07142857142857142

$ws.Range("H18").Value = 0.1
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.5
$ws.Range("O18").Value = 0.2
$ws.Range("S18").Value = 0.1

$ws.Range("F19").Value = 0.03225806451612903
$ws.Range("H19").Value = 0.2419354838709677
$ws.Range("I19").Value = 0.07258064516129033
$ws.Range("J19").Value = 0.3145161290322581
$ws.Range("K19").Value = 0.1048387096774194
$ws.Range("M19").Value = 0.03225806451612903
$ws.Range("S19").Value = 0.1370967741935484
